$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: force a freshly-created cell in column B (which has an ambiguous
# <col> style definition in this workbook) to pick up the "normal / wrap
# text / top aligned" look used throughout column B, instead of picking up
# column A's bold style. We do this by copying the formatting only (not the
# value) from a column-B cell whose style never changes in this script.
# ---------------------------------------------------------------------------
function Set-ColumnBLook($range) {
    $ws.Range("B9").Copy()
    $range.PasteSpecial(-4122)
}

# Row 10 (Objetivos:) - B/C text replaced
$ws.Range("B10").Value = "471420 - Carlos Antonio Reis Pereira Baptista"
$ws.Range("C10").Value = "471420 - Carlos Antonio Reis Pereira Baptista"

# Row 13 used to hold the first "Docentes responsaveis" entry in B/C with no
# label in A; it becomes the "Programa resumido:" row, re-using the
# "Ativacao:" date-look-alike text value in B/C. The value "01/01/2022" is
# stored as text (not a real date) in the source row B8/C8, so copy it
# across instead of retyping it (retyping would make Excel reinterpret it
# as a date serial number).
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B8").Copy()
$ws.Range("B13").PasteSpecial(-4163)
$ws.Range("C8").Copy()
$ws.Range("C13").PasteSpecial(-4163)
$ws.Rows.Item(13).RowHeight = 60

# Row 14 used to hold the second "Docentes responsaveis" entry; it becomes
# "Short syllabus:" with no B/C content.
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Clear()
$ws.Range("C14").Clear()
$ws.Rows.Item(14).RowHeight = 60

# Row 15 was "Programa resumido:" / "A definir..."; becomes "Programa:" with
# the first professor's name in B/C.
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "471420 - Carlos Antonio Reis Pereira Baptista"
$ws.Range("C15").Value = "471420 - Carlos Antonio Reis Pereira Baptista"
$ws.Rows.Item(15).RowHeight = 120

# Row 16 was "Short syllabus:"; becomes "Syllabus:"
$ws.Range("A16").Value = "Syllabus:"
$ws.Rows.Item(16).RowHeight = 120

# Row 17 was "Programa:" with a long B/C paragraph; becomes "Avaliacao:"
# with no B/C content, and loses its custom row height entirely.
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").Clear()
$ws.Range("C17").Clear()
$ws.Rows.Item(17).AutoFit()

# Row 18 was "Syllabus:" with no B/C; becomes "Metodo:" with the second
# professor's name in B/C (fresh cells, so fix up column B's look first).
$ws.Range("A18").Value = "Método:"
Set-ColumnBLook $ws.Range("B18")
$ws.Range("B18").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C18").Value = "519033 - Carlos Yujiro Shigue"
$ws.Rows.Item(18).RowHeight = 60

# Row 19 was "Avaliacao:" with no B/C; becomes "Criterio:" with the
# evaluation paragraph in B/C (fresh cells).
$ws.Range("A19").Value = "Critério:"
Set-ColumnBLook $ws.Range("B19")
$ws.Range("B19").Value = "Este curso deverá conter avaliações escritas e desenvolvimento de Estudo de Casos ou Projetos na área de Engenharia de Materiais. Sendo necessário aplicar pelo menos dois tipos de avaliações diferentes."
$ws.Range("C19").Value = "Este curso deverá conter avaliações escritas e desenvolvimento de Estudo de Casos ou Projetos na área de Engenharia de Materiais. Sendo necessário aplicar pelo menos dois tipos de avaliações diferentes."
$ws.Rows.Item(19).RowHeight = 60

# Row 20 was "Metodo:" with the evaluation paragraph; becomes "Norma de
# recuperacao:" with the M=(A1+A2)/2 text. Height stays 60.
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "A média do semestre será computada com base na relação: M=(A1+A2)/2"
$ws.Range("C20").Value = "A média do semestre será computada com base na relação: M=(A1+A2)/2"
$ws.Rows.Item(20).RowHeight = 60

# Row 21 was "Criterio:" with the M=(A1+A2)/2 text; becomes "Bibliografia:"
# with "Nao cabe recuperacao." text, and grows to 120 height.
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Não cabe recuperação."
$ws.Range("C21").Value = "Não cabe recuperação."
$ws.Rows.Item(21).RowHeight = 120

# Rows 22 and 23 (old "Norma de recuperacao:" and "Bibliografia:" rows) are
# dropped entirely, shrinking the sheet from A1:C23 down to A1:C21.
$ws.Rows("22:23").Delete()
